$d = $word.ActiveDocument

$map = [ordered]@{
  "546×4=" = "461×6="
  "807×3=" = "756×6="
  "607×8=" = "868×6="
  "783×2=" = "796×7="
  "429×5=" = "152×5="
  "591×4=" = "826×4="
  "361×3=" = "993×9="
  "458×2=" = "415×6="
  "351×9=" = "762×4="
  "220×9=" = "433×9="
  "567×5=" = "920×2="
  "662×3=" = "138×9="
  "293×2=" = "837×5="
  "400×8=" = "961×7="
  "954×9=" = "523×6="
  "978×3=" = "798×4="
  "309×3=" = "595×7="
  "921×8=" = "939×3="
  "519×5=" = "739×3="
  "153×7=" = "727×5="
  "179×4=" = "493×8="
  "564×6=" = "336×3="
  "825×9=" = "689×6="
  "489×6=" = "981×9="
  "611×3=" = "944×7="
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
